$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Carlos Pascal (new contact) - email entered before the name cells
$ws.Cells.Item(2, 3).Value = "carlos@carlos.com.br"
$ws.Cells.Item(2, 1).Value = "Carlos"
$ws.Cells.Item(2, 2).Value = "Pascal"
$ws.Cells.Item(2, 4).Value = "Junior"

# Row 3: Stenio Rapchan (existing first name kept, last name/email updated)
$ws.Cells.Item(3, 1).Value = "Stenio"
$ws.Cells.Item(3, 2).Value = "Rapchan"
$ws.Cells.Item(3, 3).Value = "stenio.rapchan@gmail.com"
$ws.Cells.Item(3, 4).Value = "All"

# Row 4: Odelinda JoJo (new contact) - email entered before the name cells
$ws.Cells.Item(4, 3).Value = "odelinda3848@uorak.com"
$ws.Cells.Item(4, 1).Value = "Odelinda"
$ws.Cells.Item(4, 2).Value = "JoJo"
$ws.Cells.Item(4, 4).Value = "Director"

# Row 5: Milton Mestre (names kept, email updated)
$ws.Cells.Item(5, 1).Value = "Milton"
$ws.Cells.Item(5, 2).Value = "Mestre"
$ws.Cells.Item(5, 3).Value = "miltinho@gmail.com"
$ws.Cells.Item(5, 4).Value = "Expert"

# Row 6: Ravi Messias (names kept, email updated)
$ws.Cells.Item(6, 1).Value = "Ravi"
$ws.Cells.Item(6, 2).Value = "Messias"
$ws.Cells.Item(6, 3).Value = "ravi@hotmail.com"
$ws.Cells.Item(6, 4).Value = "Intermediate"

# Row 7: Rodrigo Rodrigues (new contact)
$ws.Cells.Item(7, 1).Value = "Rodrigo"
$ws.Cells.Item(7, 2).Value = "Rodrigues"
$ws.Cells.Item(7, 3).Value = "roro@uol.com.br"
$ws.Cells.Item(7, 4).Value = "Junior"

# Row 8: Reginaldo Maromba (new row)
$ws.Cells.Item(8, 1).Value = "Reginaldo"
$ws.Cells.Item(8, 2).Value = "Maromba"
$ws.Cells.Item(8, 3).Value = "regis100@hotmail.com"
$ws.Cells.Item(8, 4).Value = "Intermediate"

# Row 9: Guilhermino Novais (new row)
$ws.Cells.Item(9, 1).Value = "Guilhermino"
$ws.Cells.Item(9, 2).Value = "Novais"
$ws.Cells.Item(9, 3).Value = "guinova@gmail.com"
$ws.Cells.Item(9, 4).Value = "Expert"

# Columns A and B re-fit to the new (longer) content.
$ws.Columns.Item(1).ColumnWidth = 11.25
$ws.Columns.Item(2).ColumnWidth = 9.25

$ws.Range("C24").Select()
